$wb = $excel.ActiveWorkbook
$dataSheet = $wb.Worksheets.Item(1)

# --- Update the "panel_query_time" values (column F) on the "data" sheet ---
$newTimes = @(
  "2021-10-05 14:21:51.332521",
  "2021-10-05 14:21:51.332529",
  "2021-10-05 14:21:51.332533",
  "2021-10-05 14:21:51.332535",
  "2021-10-05 14:21:51.332538",
  "2021-10-05 14:21:51.332541",
  "2021-10-05 14:21:51.332544",
  "2021-10-05 14:21:51.332546",
  "2021-10-05 14:21:51.332549",
  "2021-10-05 14:21:51.332551",
  "2021-10-05 14:21:51.332554",
  "2021-10-05 14:21:51.332556",
  "2021-10-05 14:21:51.332559",
  "2021-10-05 14:21:51.332561",
  "2021-10-05 14:21:51.332564",
  "2021-10-05 14:21:51.332567",
  "2021-10-05 14:21:51.332570",
  "2021-10-05 14:21:51.332572",
  "2021-10-05 14:21:51.332575",
  "2021-10-05 14:21:51.332577",
  "2021-10-05 14:21:51.332580",
  "2021-10-05 14:21:51.332582",
  "2021-10-05 14:21:51.332585",
  "2021-10-05 14:21:51.332587",
  "2021-10-05 14:21:51.332590",
  "2021-10-05 14:21:51.332593",
  "2021-10-05 14:21:51.332596",
  "2021-10-05 14:21:51.332598",
  "2021-10-05 14:21:51.332601",
  "2021-10-05 14:21:51.332603",
  "2021-10-05 14:21:51.332606",
  "2021-10-05 14:21:51.332608",
  "2021-10-05 14:21:51.332611",
  "2021-10-05 14:21:51.332614",
  "2021-10-05 14:21:51.332616",
  "2021-10-05 14:21:51.332619",
  "2021-10-05 14:21:51.332621",
  "2021-10-05 14:21:51.332624",
  "2021-10-05 14:21:51.332626",
  "2021-10-05 14:21:51.332629",
  "2021-10-05 14:21:51.332632",
  "2021-10-05 14:21:51.332635",
  "2021-10-05 14:21:51.332637",
  "2021-10-05 14:21:51.332640",
  "2021-10-05 14:21:51.332642",
  "2021-10-05 14:21:51.332645",
  "2021-10-05 14:21:51.332647",
  "2021-10-05 14:21:51.332650",
  "2021-10-05 14:21:51.332652",
  "2021-10-05 14:21:51.332655",
  "2021-10-05 14:21:51.332657"
)

for ($i = 0; $i -lt $newTimes.Count; $i++) {
  $row = $i + 2
  $dataSheet.Cells.Item($row, 6).Value = $newTimes[$i]
}

# --- Add the new "metadata" sheet, placed right after "data" ---
$meta = $wb.Worksheets.Add($null, $dataSheet)
$meta.Name = "metadata"

# Copy the header-row / index-column formatting from "data" before writing values
# so the new cells pick up the same bold/border/centered style.
$dataSheet.Range("B1:F1").Copy()
$meta.Range("B1:G1").PasteSpecial(-4122)
$dataSheet.Range("A2").Copy()
$meta.Range("A2").PasteSpecial(-4122)

# Header row
$meta.Cells.Item(1, 2).Value = "data_name"
$meta.Cells.Item(1, 3).Value = "data_id"
$meta.Cells.Item(1, 4).Value = "data_version"
$meta.Cells.Item(1, 5).Value = "data_version_created"
$meta.Cells.Item(1, 6).Value = "panel_query_time"
$meta.Cells.Item(1, 7).Value = "panel_get_request"

# Data row
$meta.Cells.Item(2, 1).Value = 0
$meta.Cells.Item(2, 2).Value = "Non-syndromic familial congenital anorectal malformations"
$meta.Cells.Item(2, 3).Value = 253

# "data_version" (1.7) must be stored as text, not a number - force text entry
# via a temporary "@" number format, then restore the cell to the default
# (unstyled) look so it matches the plain/unstyled data cells around it.
$meta.Cells.Item(2, 4).NumberFormat = "@"
$meta.Cells.Item(2, 4).Value = "1.7"
$dataSheet.Range("A1").Copy()
$meta.Range("D2").PasteSpecial(-4122)

$meta.Cells.Item(2, 5).Value = "2021-10-04T13:41:22.431587Z"
$meta.Cells.Item(2, 6).Value = "2021-10-05 14:21:51.328751"
$meta.Cells.Item(2, 7).Value = "https://panelapp.genomicsengland.co.uk/api/v1/panels/253/?format=json"

$dataSheet.Activate()
